$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unneeded rows (old rows 7-11) and column F (Nilai Defuzzifikasi)
$ws.Rows("7:11").Delete()
$ws.Columns("F:F").Delete()

# Update header row (row 1)
$ws.Range("B1").Value = "PELAYANAN"
$ws.Range("C1").Value = "HARGA"
$ws.Range("D1").Value = "SKOR_KELAYAKAN"
$ws.Range("E1").Value = "KETERANGAN"

# Re-populate data rows 2-6 with the new (re-ranked) values
$ws.Range("A2").Value = 79
$ws.Range("B2").Value = 92
$ws.Range("C2").Value = 22360
$ws.Range("D2").Value = 92
$ws.Range("E2").Value = "Sangat Layak"

$ws.Range("A3").Value = 80
$ws.Range("B3").Value = 89
$ws.Range("C3").Value = 22012
$ws.Range("D3").Value = 89
$ws.Range("E3").Value = "Sangat Layak"

$ws.Range("A4").Value = 25
$ws.Range("B4").Value = 94
$ws.Range("C4").Value = 34513
$ws.Range("D4").Value = 87.50000000000001
$ws.Range("E4").Value = "Sangat Layak"

$ws.Range("A5").Value = 78
$ws.Range("B5").Value = 86
$ws.Range("C5").Value = 27315
$ws.Range("D5").Value = 86
$ws.Range("E5").Value = "Sangat Layak"

$ws.Range("A6").Value = 86
$ws.Range("B6").Value = 84
$ws.Range("C6").Value = 29811
$ws.Range("D6").Value = 85.2400728182956
$ws.Range("E6").Value = "Sangat Layak"
